$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) and volume (E) columns keep their text formatting so that
# values such as "323.02" or "  +0.65%  " are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.777.80"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "1.857.74"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "323.02"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "1.032"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.4405"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("D9").Value = "0.07446"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").Value = "0.8878"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("D11").Value = "21.61"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "1.863.62"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "5.533"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "6.754"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "0.07210"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "86.10"
$ws.Range("E16").Value = "  +4.11%  "
$ws.Range("D17").Value = "1.039"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "0.000009121"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").Value = "1.033"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "15.58"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "27.817.79"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "5.306"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("D24").Value = "2.091.25"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").Value = "2.078"
$ws.Range("E25").Value = "  +6.77%  "
$ws.Range("D26").Value = "159.10"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "18.77"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "5.362"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("D29").Value = "1.996"
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("D30").Value = "118.96"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("D31").Value = "0.09121"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").Value = "1.216"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").Value = "0.7746"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").Value = "3.042"
$ws.Range("E34").Value = "  +5.46%  "
$ws.Range("D35").Value = "4.605"
$ws.Range("E35").Value = "  +1.58%  "
$ws.Range("D36").Value = "1.034"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "1.157"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "0.01989"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").Value = "0.05326"
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("D40").Value = "2.863"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("D41").Value = "0.5221"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("D42").Value = "6.963"
$ws.Range("E42").Value = "  +3.10%  "
$ws.Range("D43").Value = "0.1680"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").Value = "8.815"
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "110.82"
$ws.Range("E45").Value = "  +1.62%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "10.79"
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("D47").Value = "1.035"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").Value = "0.06564"
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("D49").Value = "1.718"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "0.4740"
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("D51").Value = "1.892"
$ws.Range("E51").Value = "  +0.55%  "
